$d = $word.ActiveDocument

$startRng = $d.Content
$startRng.Find.Execute("PROTÓTIPOS") | Out-Null
$delStart = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("REFERÊNCIAS") | Out-Null
$delEnd = $endRng.Start

Write-Host ("delStart: " + $delStart + " delEnd: " + $delEnd)

$delRange = $d.Range($delStart, $delEnd)
Write-Host ("text to delete len: " + $delRange.Text.Length)
$delRange.Delete()
Write-Host ("Paragraphs after: " + $d.Paragraphs.Count)
